# Generate Report for Handback
# Updates the localization-status workbook: the de-de/zh-cn handback has
# landed and is now in sync with en-US, so refresh the Status, the
# "Latest Handback DateTime" stamps, and clear the stale "handback file is
# not the latest" Error Detail message for both locale sheets. Column
# widths on the affected columns are re-fit to the (now different) cell
# content, same as Excel does after a content refresh.

$wb = $excel.ActiveWorkbook

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("K2").Value = "2016-09-06 17:26:34"
$wsZh.Range("P2").Value = ""

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("K2").Value = "2016-09-06 17:26:42"
$wsDe.Range("P2").Value = ""

# ---- re-fit the columns whose content changed length ----
# Status (now the longer "Handed back: in sync with en-US") needs a wider
# column; Error Detail (now blank) needs a narrower one. The "zh-cn"/"de-de"
# Status column mirrors onto the Overview sheet's per-locale status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(16).ColumnWidth = 12.8

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(16).ColumnWidth = 12.8
